$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# DCT_RIGHTSHOLDER lives in column Q. Make the column mandatory:
# highlight the header + example rows with the same "required" yellow
# fill used by the other mandatory columns, and remove the sample
# values ("jan", "eko", "linda") that used to populate the example rows.
$ws.Range("Q2").ClearContents()
$ws.Range("Q3").ClearContents()
$ws.Range("Q4").ClearContents()

$ws.Range("Q1:Q4").Interior.Color = 65535

# Reflect where the user was working when they made the column mandatory.
$ws.Range("Q1").Select() | Out-Null
